$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header cells
$ws.Range("H1").Value = "Status of cases"
$ws.Range("I1").Value = "31-day average of new deaths per million on end_date"
$ws.Range("J1").Value = "61-day average of new deaths per million on end_date"

# Populate new I (31-day avg deaths) and J (61-day avg deaths) columns
$ws.Range("I2").Value = 0.0455483870967742
$ws.Range("J2").Value = 0.04672131147540984
$ws.Range("I3").Value = 0.1588387096774194
$ws.Range("J3").Value = 0.1779016393442623
$ws.Range("I4").Value = 0.09125806451612903
$ws.Range("J4").Value = 0.07027868852459017
$ws.Range("I5").Value = 0.1663548387096774
$ws.Range("J5").Value = 0.1973934426229508
$ws.Range("I6").Value = 0.002645161290322581
$ws.Range("J6").Value = 0.002688524590163935
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("I8").Value = 2.799451612903226
$ws.Range("J8").Value = 6.184934426229509
$ws.Range("I9").Value = 0.01393548387096774
$ws.Range("J9").Value = 0.008655737704918034
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("I12").Value = 0.009709677419354839
$ws.Range("J12").Value = 0.01357377049180328
$ws.Range("I13").Value = 1.972806451612903
$ws.Range("J13").Value = 1.326983606557377
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0.003393442622950819
$ws.Range("I15").Value = 0.02164516129032258
$ws.Range("J15").Value = 0.016
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("I17").Value = 0.02335483870967742
$ws.Range("J17").Value = 0.09504918032786885
$ws.Range("I18").Value = 0.001225806451612903
$ws.Range("J18").Value = 0.008098360655737704
$ws.Range("I19").Value = 0.005709677419354839
$ws.Range("J19").Value = 0.007409836065573771
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0.03318032786885246
$ws.Range("I21").Value = 0.1345161290322581
$ws.Range("J21").Value = 0.151672131147541
$ws.Range("I22").Value = 0.5968709677419355
$ws.Range("J22").Value = 0.7582786885245901
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("I24").Value = 0.08503225806451613
$ws.Range("J24").Value = 0.1163606557377049
$ws.Range("I25").Value = 0.2003225806451613
$ws.Range("J25").Value = 0.4341803278688525
$ws.Range("I26").Value = 0.01654838709677419
$ws.Range("J26").Value = 0.04162295081967214
$ws.Range("I27").Value = 0.01716129032258065
$ws.Range("J27").Value = 0.02491803278688524
$ws.Range("I28").Value = 0.01638709677419355
$ws.Range("J28").Value = 0.05831147540983606
$ws.Range("I29").Value = 0.03393548387096774
$ws.Range("J29").Value = 0.05031147540983606
$ws.Range("I30").Value = 1.51441935483871
$ws.Range("J30").Value = 1.630327868852459
$ws.Range("I31").Value = 0.7336451612903225
$ws.Range("J31").Value = 0.738655737704918
$ws.Range("I32").Value = 0.109741935483871
$ws.Range("J32").Value = 0.1015245901639344
$ws.Range("I33").Value = 0.197741935483871
$ws.Range("J33").Value = -0.9850163934426229
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("I35").Value = 0.1355806451612903
$ws.Range("J35").Value = 0.1301311475409836
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("I37").Value = 0.02787096774193548
$ws.Range("J37").Value = 0.04190163934426229
$ws.Range("I38").Value = 0.005032258064516129
$ws.Range("J38").Value = 0.02052459016393442
$ws.Range("I39").Value = 0.00632258064516129
$ws.Range("J39").Value = 0.005622950819672131
$ws.Range("I40").Value = 0.01387096774193548
$ws.Range("J40").Value = 0.02114754098360656
$ws.Range("I41").Value = 3.198612903225806
$ws.Range("J41").Value = 2.580508196721312
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("I43").Value = 1.012129032258065
$ws.Range("J43").Value = 0.9504918032786885
$ws.Range("I44").Value = 0.0392258064516129
$ws.Range("J44").Value = 0.02832786885245902
$ws.Range("I45").Value = 0.4345806451612903
$ws.Range("J45").Value = 0.228672131147541
$ws.Range("I46").Value = 0.3698064516129032
$ws.Range("J46").Value = 0.3331475409836066
$ws.Range("I47").Value = 0.04870967741935484
$ws.Range("J47").Value = 0.06436065573770491
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("I49").Value = 0.005290322580645161
$ws.Range("J49").Value = 0.01139344262295082
$ws.Range("I50").Value = 0.03390322580645162
$ws.Range("J50").Value = 0.03367213114754098
$ws.Range("I51").Value = 1.378322580645161
$ws.Range("J51").Value = 1.13744262295082
$ws.Range("I52").Value = 0.003612903225806452
$ws.Range("J52").Value = 0.007344262295081967
$ws.Range("I53").Value = 0.5496451612903226
$ws.Range("J53").Value = 0.5830163934426229
$ws.Range("I54").Value = 0.02732258064516129
$ws.Range("J54").Value = 0.03155737704918032
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("I56").Value = 0.03670967741935484
$ws.Range("J56").Value = 0.06363934426229508
$ws.Range("I57").Value = 0.004032258064516129
$ws.Range("J57").Value = 0.00819672131147541
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("I59").Value = 0.002032258064516129
$ws.Range("J59").Value = 0.00619672131147541
$ws.Range("I60").Value = 0.01725806451612903
$ws.Range("J60").Value = 0.01168852459016393
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0.001540983606557377
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0.01419672131147541
$ws.Range("I63").Value = 0.3892903225806452
$ws.Range("J63").Value = 0.6500327868852459
$ws.Range("I64").Value = 0.1490645161290323
$ws.Range("J64").Value = 0.1665737704918033
$ws.Range("I65").Value = 0.02709677419354839
$ws.Range("J65").Value = 0.02754098360655738
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("I68").Value = 0.04290322580645161
$ws.Range("J68").Value = 0.04757377049180328
$ws.Range("I69").Value = 1.078161290322581
$ws.Range("J69").Value = 0.6352950819672131
$ws.Range("I70").Value = 0.029
$ws.Range("J70").Value = 0.02985245901639344
$ws.Range("I71").Value = 1.54241935483871
$ws.Range("J71").Value = 1.221295081967213
$ws.Range("I72").Value = 0.1148064516129032
$ws.Range("J72").Value = 0.1425245901639344
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0.002180327868852459
$ws.Range("I74").Value = 0.01532258064516129
$ws.Range("J74").Value = 0.03757377049180328
$ws.Range("I75").Value = 0.0455483870967742
$ws.Range("J75").Value = 0.07662295081967213
$ws.Range("I76").Value = 0.01512903225806452
$ws.Range("J76").Value = 0.1113934426229508
